# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# handback/report generation has completed for both the zh-cn and de-de
# languages:
#   - Overview status text changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" (affects all cells showing that text).
#   - The (previously blank) "Latest Target File" / "Latest Handback File"
#     columns on the per-locale sheets are populated with hyperlinked file
#     names and handback datetimes.
#   - Several columns are widened to better fit the new content.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Overview sheet: update status text + widen the status columns
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------
# 2. zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# The Status column for both rows also reads "Ready for handoff" -> update
# it to stay in sync with the shared text change above.
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

# Widen Status / Latest Target File / Latest Handback File columns
$zhcn.Columns.Item(3).ColumnWidth = 29.1
$zhcn.Columns.Item(9).ColumnWidth = 39.1
$zhcn.Columns.Item(10).ColumnWidth = 39.1

# Row 2 (a5eaf1b1...): populate Latest Target File (I) with a hyperlink to
# the source .md file, Latest Handback File (J) with the generated xlf, and
# update the Latest Handback DateTime (K).
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c2c0219b9058da9d087bad0c42b6fcbb7c174c0/e2e/a5eaf1b1-894e-473a-94fa-a2367a38b7ac.md", "", "", "a5eaf1b1-894e-473a-94fa-a2367a38b7ac.md")
$zhcn.Range("J2").Value = "a5eaf1b1-894e-473a-94fa-a2367a38b7ac.7ac40d134b778e577f94f8c45fb09a88968b53d8.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-06 06:46:22"

# Row 3 (b8174788...): same treatment.
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c2c0219b9058da9d087bad0c42b6fcbb7c174c0/e2e/b8174788-9a90-4227-8136-0f93962b431e.md", "", "", "b8174788-9a90-4227-8136-0f93962b431e.md")
$zhcn.Range("J3").Value = "b8174788-9a90-4227-8136-0f93962b431e.b81b39c77959a143b41bf4b3840b5e702d70060a.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-06 06:46:22"

# ---------------------------------------------------------------------
# 3. de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

# The Status column for both rows also reads "Ready for handoff" -> update
# it to stay in sync with the shared text change above.
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

# Widen Status / Latest Target File / Latest Handback File columns
$dede.Columns.Item(3).ColumnWidth = 29.1
$dede.Columns.Item(9).ColumnWidth = 39.1
$dede.Columns.Item(10).ColumnWidth = 39.1

# Row 2 (a5eaf1b1...)
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c2c0219b9058da9d087bad0c42b6fcbb7c174c0/e2e/a5eaf1b1-894e-473a-94fa-a2367a38b7ac.md", "", "", "a5eaf1b1-894e-473a-94fa-a2367a38b7ac.md")
$dede.Range("J2").Value = "a5eaf1b1-894e-473a-94fa-a2367a38b7ac.7ac40d134b778e577f94f8c45fb09a88968b53d8.de-de.xlf"
$dede.Range("K2").Value = "2016-09-06 06:46:41"

# Row 3 (b8174788...)
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c2c0219b9058da9d087bad0c42b6fcbb7c174c0/e2e/b8174788-9a90-4227-8136-0f93962b431e.md", "", "", "b8174788-9a90-4227-8136-0f93962b431e.md")
$dede.Range("J3").Value = "b8174788-9a90-4227-8136-0f93962b431e.b81b39c77959a143b41bf4b3840b5e702d70060a.de-de.xlf"
$dede.Range("K3").Value = "2016-09-06 06:46:41"

$wb.Save()
